# Update "想去人数" (want-to-go count) values in column F on both the
# "展览" sheet and the "全部类型" sheet for the affected events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): rows 2-6
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value = 1045
$wsExh.Range("F3").Value = 263
$wsExh.Range("F4").Value = 2672
$wsExh.Range("F6").Value = 570

# Sheet "全部类型" (All types): rows 4-8 (same events, shifted by 2 rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1045
$wsAll.Range("F5").Value = 263
$wsAll.Range("F6").Value = 2672
$wsAll.Range("F8").Value = 570
